$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Three table shapes (slides 14, 15, 16) switch from the default table
#    style {81DB735D-0F7E-4BCE-B422-E3DB8D7ED1C1} to the built-in style
#    {BEDED662-EAC5-4EA3-933E-899FE628DEB7}.
# ---------------------------------------------------------------------------
$targetStyle = "{BEDED662-EAC5-4EA3-933E-899FE628DEB7}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyle)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) The deck's colour theme is swapped from the "Integral" (Red Violet)
#    palette to the default "Office Theme" palette. Re-point every theme
#    colour slot at the Office values via the live ThemeColorScheme (this
#    rewrites ppt/theme/theme1.xml, which is the part the slide master
#    actually renders with).
# ---------------------------------------------------------------------------
$officeColors = @{
    1  = 0x000000  # dk1
    2  = 0xFFFFFF  # lt1
    3  = 0x6A5444  # dk2       (44546A, stored BGR)
    4  = 0xE6E6E7  # lt2       (E7E6E6, stored BGR)
    5  = 0xD59B5B  # accent1   (5B9BD5, stored BGR)
    6  = 0x317DED  # accent2   (ED7D31, stored BGR)
    7  = 0xA5A5A5  # accent3   (A5A5A5, stored BGR)
    8  = 0x00C0FF  # accent4   (FFC000, stored BGR)
    9  = 0xC47244  # accent5   (4472C4, stored BGR)
    10 = 0x47AD70  # accent6   (70AD47, stored BGR)
    11 = 0xC16305  # hlink     (0563C1, stored BGR)
    12 = 0x724F95  # folHlink  (954F72, stored BGR)
}

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i]
}
